$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'ski knee',
    'ski knee support',
    'ski pant',
    'ski pant patch',
    'ski pants patch',
    'ski pants women',
    'ski pants women black',
    'ski pants women short length',
    'ski pants women xs short',
    'ski skins',
    'ski tops',
    'ski tops women',
    'ski training',
    'ski training shorts',
    'ski wear',
    'ski wear women',
    'skiing clothing women',
    'skiing gear',
    'skiing knee brace',
    'skiing leggings',
    'skiing pants',
    'skiing pants for women',
    'skiing suit women',
    'skiing tights',
    'skims cozy loungewear',
    'skims loungewear',
    'skin color leggings',
    'skin color leggings for men',
    'skin color leggings for women',
    'skin color tights',
    'skin color tights for women',
    'skin compression',
    'skin compression men',
    'skin compression tights',
    'skin compression women',
    'skin pants women',
    'skin saver international',
    'skin stockings for women',
    'skin tight clothes for women',
    'skin tight clothes women',
    'skin tight dress for women',
    'skin tight leggings',
    'skin tight tops for women',
    'skin tights',
    'skin tights for women',
    'skin to fabric tape',
    'skin to skin top',
    'skin tone compression tights',
    'skin tone leggings for women',
    'skins a400',
    'skins a400 compression',
    'skins a400 long tights',
    'skins a400 women',
    'skins compression',
    'skins compression men',
    'skins compression men tights',
    'skins compression tights men',
    'skins compression tights women',
    'skins compression top',
    'skins compression wear',
    'skins compression women',
    'skins dnamic compression long tights',
    'skins dnamic women''s compression tights',
    'skins dynamic long mens tights',
    'skins dynamic tights',
    'skins half tights',
    'skins leggings',
    'skins running tights men',
    'skins tights',
    'skins tights women',
    'skins women compression',
    'skirt and leggings',
    'skirt and leggings for women',
    'skirt and pants leggings',
    'skirt capri',
    'skirt capri leggings',
    'skirt hike',
    'skirt legging',
    'skirt leggings for women',
    'skirt leggings for women sport',
    'skirt leggins',
    'skirt with capri leggings',
    'skirt with leggings',
    'skirt with leggings women',
    'skirts over leggings',
    'skirts with leggings for women',
    'skis womens',
    'sleep knee support',
    'sleep lady',
    'sleep leggings for women',
    'sleep tights women',
    'sleep wear womens',
    'slick leggings for women',
    'slide pants',
    'slimmer thighs',
    'slip joint',
    'small black',
    'small clothes',
    'small generator',
    'small knee support'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
